$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header cell B1: "Typenprüfung" -> "Typenprüfung &\nKonvertierung" (now wrapped, bold, 2 lines)
$ws.Range("B1").Value = "Typenprüfung &`nKonvertierung"
$ws.Range("B1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 30.75

# Conversion-function column (B) updates: the generic "bin2hex" placeholder is
# replaced by type-specific conversion routines.
$ws.Range("B2").Value = "bit->hex"        # BIT
$ws.Range("B17").Value = "bin->hex"       # BINARY
$ws.Range("B18").Value = "base64->hex"    # VARBINARY
$ws.Range("B19").Value = "base64->hex"    # LONGVARBINARY
$ws.Range("B26").Value = "base64->hex"    # BLOB

# Restore the cursor/selection position as recorded in the saved workbook.
[void]$ws.Range("J17").Select()
